$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = -0.32703571387048669
$ws.Range("B1").Value = 0.32666121046958807
$ws.Range("A2").Value = -0.2051753426779328
$ws.Range("B2").Value = 0.20450361354225954
$ws.Range("A3").Value = -0.15479719613630039
$ws.Range("B3").Value = 0.15444368530624608
$ws.Range("A4").Value = -0.14644368536363572
$ws.Range("B4").Value = 0.14566787571483353
$ws.Range("A5").Value = -0.1426678757511084
$ws.Range("B5").Value = 0.13998682748191182
$ws.Range("A6").Value = -0.040738420651523555
$ws.Range("B6").Value = 0.040314636685835126
$ws.Range("A7").Value = -0.030314636772692527
$ws.Range("B7").Value = 0.030216873686327084
$ws.Range("A8").Value = -0.020216873775021682
$ws.Range("B8").Value = 0.020058221011908905
$ws.Range("A9").Value = -0.018058221061724833
$ws.Range("B9").Value = 0.017933212265003107
$ws.Range("A10").Value = -0.01593321231674949
$ws.Range("B10").Value = 0.01592608728119238
$ws.Range("A11").Value = -0.012926087338200887
$ws.Range("B11").Value = 0.012914415902192111
$ws.Range("A12").Value = -0.0094144159621416534
$ws.Range("B12").Value = 0.0093368099443762453
$ws.Range("A13").Value = -0.0058368100057935024
$ws.Range("B13").Value = 0.0058077073647471167
$ws.Range("A14").Value = 0.0021922925503936597
$ws.Range("B14").Value = -0.0021970452353414771
$ws.Range("A15").Value = -0.0080538607823887887
$ws.Range("B15").Value = 0.0080350391820775968
$ws.Range("A16").Value = -0.006035039236616857
$ws.Range("B16").Value = 0.0060035683054908695
$ws.Range("A17").Value = -0.0040035683608721229
$ws.Range("B17").Value = 0.0039999999341215897
$ws.Range("A18").Value = -0.088631607072812812
$ws.Range("B18").Value = 0.088499942654042485
$ws.Range("A19").Value = -0.084499942677639162
$ws.Range("B19").Value = 0.083505602082860175
$ws.Range("A20").Value = -0.064594089427773227
$ws.Range("B20").Value = 0.06430614269502577
$ws.Range("A21").Value = -0.0040058490143737302
$ws.Range("B21").Value = 0.0039999999622608584
$ws.Range("A22").Value = -0.045706417438184488
$ws.Range("B22").Value = 0.045494741540949235
$ws.Range("A23").Value = -0.040494741579833793
$ws.Range("B23").Value = 0.040098178697327569
$ws.Range("A24").Value = -0.020098178816772005
$ws.Range("B24").Value = 0.019999999879015462
$ws.Range("A25").Value = -0.097248406875225513
$ws.Range("B25").Value = 0.097123882855983723
$ws.Range("A26").Value = -0.094623882904864232
$ws.Range("B26").Value = 0.094463479747988188
$ws.Range("A27").Value = -0.091963479799202386
$ws.Range("B27").Value = 0.091013154573710953
$ws.Range("A28").Value = -0.089013154633355462
$ws.Range("B28").Value = 0.0883625315878378
$ws.Range("A29").Value = -0.081362531680928107
$ws.Range("B29").Value = 0.081173891405348364
$ws.Range("A30").Value = -0.021173891774512565
$ws.Range("B30").Value = 0.021023871987341458
$ws.Range("A31").Value = -0.014023872086658784
$ws.Range("B31").Value = 0.014001320725942179
$ws.Range("A32").Value = -0.0040013208410378809
$ws.Range("B32").Value = 0.0039999999155817534
